$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing mobile numbers (column C, rows 2-6) ---
$ws.Range("C2").Value = 9928788888
$ws.Range("C3").Value = 8888599252
$ws.Range("C4").Value = 6665299999
$ws.Range("C5").Value = 1355555000
$ws.Range("C6").Value = 9444440000

# --- Append two new data rows (7 & 8), copying row 6's formatting ---
$ws.Rows(6).Copy()
$ws.Rows(7).Insert()
$ws.Rows(6).Copy()
$ws.Rows(8).Insert()

# Row 7: VP Traders / Vishal patel
$ws.Range("A7").Value = "VP Traders"
$ws.Range("B7").Value = "Vishal patel"
$ws.Range("C7").Value = 4878487500
$ws.Range("D7").Value = "maninagar"
$ws.Range("E7").Value = "Ahmedabad"
$ws.Range("F7").Value = "Active"

# Row 8: PV Traders / Pooja patel
$ws.Range("A8").Value = "PV Traders"
$ws.Range("B8").Value = "Pooja patel"
$ws.Range("C8").Value = 4878487501
$ws.Range("D8").Value = "maninagar"
$ws.Range("E8").Value = "Ahmedabad"
$ws.Range("F8").Value = "Active"

# --- Restore selection to match target (active cell E3) ---
$null = $ws.Range("E3").Select()
